$d = $word.ActiveDocument
$vt = [char]11
$brk = "$vt$vt"

function Set-ParaText($idx, $text) {
    $p = $d.Paragraphs.Item($idx)
    $s = $p.Range.Start
    $e = $p.Range.End
    $rng = $d.Range($s, $e)
    $rng.Text = $text
}

# --- Paragraph 1: Title ---
Set-ParaText 1 "The Profound Interplay Between Chemistry and Biology: Unveiling the Symphonies of Life"

# --- Paragraph 2: Author name ---
Set-ParaText 2 "Dr. Rebecca Watts"

# --- Paragraph 3: Email ---
Set-ParaText 3 "rebecca.watts@school.edu"

# --- Paragraph 5: Body text (big paragraph) ---
$body = "The intricate tapestry of life is a symphony of chemical reactions. From the smallest cellular processes to the vast ecosystems that shape our planet, chemistry and biology intertwine, composing a symphony of life that is mesmerizing in its complexity. Understanding the interplay between these two disciplines grants us the power to unravel the secrets of life and unlock the potential for incredible discoveries."
$body += $brk
$body += "Within the minuscule realm of cells, chemistry governs the intricate dance of molecules. It orchestrates the synthesis of proteins, the metabolism of nutrients, and the replication of DNA, the very building blocks of life. Beyond the cellular level, chemistry shapes the interactions between organisms, determining their ecological niches and the intricate webs of relationships that define ecosystems."
$body += $brk
$body += "Biology and chemistry collaborate in awe-inspiring ways to maintain the delicate balance of life. Biological processes, such as photosynthesis and respiration, harness chemical energy to fuel the growth and reproduction of organisms. In turn, chemical reactions regulate biological processes, ensuring the proper function of cells, tissues, and organs. This symphony of interactions sustains the intricate equilibrium of life on Earth."
$body += $brk
$body += "Body:"
$body += $brk
$body += "The realm of chemistry unveils the secrets of matter and its transformations. Chemistry reveals the fundamental building blocks of substances, their properties, and the forces that drive their interactions. By understanding chemical principles, scientists can manipulate molecules to create new materials, medicines, and technologies that enhance our lives."
$body += $brk
$body += "Biology unveils the mysteries of living organisms, from their smallest components to the vast ecosystems they inhabit. It explores the intricate mechanisms of life, from cellular processes to complex behaviors. By understanding biological principles, scientists can gain insights into the origins of life, the evolution of species, and the delicate balance of ecosystems."
$body += $brk
$body += "When chemistry and biology unite, they unlock extraordinary possibilities. Medicinal chemistry, for instance, blends chemical principles with biological knowledge to design drugs that target specific diseases. Agricultural chemistry harnesses chemical reactions to enhance crop yields and protect plants from pests. Environmental chemistry addresses the impact of human activities on the natural world, seeking solutions to preserve and restore ecosystems."

Set-ParaText 5 $body

# --- Paragraph 7: Summary body (after "Summary" heading paragraph 6) ---
$summary = "Chemistry and biology form a seamless dance of life, intertwining their principles to orchestrate the symphony of living systems. The study of chemistry uncovers the secrets of matter and its transformations, while biology unravels the mysteries of living organisms. Together, they empower scientists to explore the depths of life, unlocking the potential for innovative discoveries, advancements in medicine, and solutions to global challenges. By delving into this profound interplay, we gain a deeper appreciation for the intricate beauty of life and the transformative power of scientific inquiry."
Set-ParaText 7 $summary

# --- Add a new empty paragraph at the very end of the document ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRng.InsertParagraphAfter()
